$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 11.29255533333333
$ws.Cells.Item(2, 8).Value = 33.877666
$ws.Cells.Item(2, 9).Value = 0.5495662219753726
$ws.Cells.Item(2, 10).Value = 0.6375557499803809
$ws.Cells.Item(2, 13).Value = 38.10639333333333
$ws.Cells.Item(2, 14).Value = 114.31918
$ws.Cells.Item(2, 15).Value = 0.3831479157160237
$ws.Cells.Item(2, 16).Value = 0.4159903984418967
$ws.Cells.Item(2, 17).Value = 430.3185552704311
$ws.Cells.Item(2, 18).Value = 3872.86699743388
$ws.Cells.Item(2, 19).Value = 0.2105651524977937
$ws.Cells.Item(2, 20).Value = 0.2652170704632609
$ws.Cells.Item(3, 7).Value = 11.29255533333333
$ws.Cells.Item(3, 8).Value = 33.877666
$ws.Cells.Item(3, 9).Value = 0.5495662219753726
$ws.Cells.Item(3, 10).Value = 0.6375557499803809
$ws.Cells.Item(3, 15).Value = 0.09199521176963764
$ws.Cells.Item(3, 16).Value = 0.09988081163714851
$ws.Cells.Item(3, 17).Value = 103.321054341448
$ws.Cells.Item(3, 18).Value = 929.8894890730321
$ws.Cells.Item(3, 19).Value = 0.05055746097206409
$ws.Cells.Item(3, 20).Value = 0.06367958577197137
$ws.Cells.Item(4, 7).Value = 11.29255533333333
$ws.Cells.Item(4, 8).Value = 33.877666
$ws.Cells.Item(4, 9).Value = 0.5495662219753726
$ws.Cells.Item(4, 10).Value = 0.6375557499803809
$ws.Cells.Item(4, 13).Value = 15.023598
$ws.Cells.Item(4, 14).Value = 45.070794
$ws.Cells.Item(4, 15).Value = 0.1510575983904562
$ws.Cells.Item(4, 16).Value = 0.1640058785774412
$ws.Cells.Item(4, 17).Value = 169.654811720756
$ws.Cells.Item(4, 18).Value = 1526.893305486804
$ws.Cells.Item(4, 19).Value = 0.08301615364811614
$ws.Cells.Item(4, 20).Value = 0.1045628909176318
$ws.Cells.Item(5, 7).Value = 11.29255533333333
$ws.Cells.Item(5, 8).Value = 33.877666
$ws.Cells.Item(5, 9).Value = 0.5495662219753726
$ws.Cells.Item(5, 10).Value = 0.6375557499803809
$ws.Cells.Item(5, 13).Value = 23.556204
$ws.Cells.Item(5, 14).Value = 47.112408
$ws.Cells.Item(5, 15).Value = 0.236850294013169
$ws.Cells.Item(5, 16).Value = 0.1714350065796238
$ws.Cells.Item(5, 17).Value = 266.0097371132881
$ws.Cells.Item(5, 18).Value = 1596.058422679728
$ws.Cells.Item(5, 19).Value = 0.1301649212545735
$ws.Cells.Item(5, 20).Value = 0.1092993741927636
$ws.Cells.Item(6, 7).Value = 11.29255533333333
$ws.Cells.Item(6, 8).Value = 33.877666
$ws.Cells.Item(6, 9).Value = 0.5495662219753726
$ws.Cells.Item(6, 10).Value = 0.6375557499803809
$ws.Cells.Item(6, 13).Value = 13.62041
$ws.Cells.Item(6, 14).Value = 40.86123000000001
$ws.Cells.Item(6, 15).Value = 0.1369489801107134
$ws.Cells.Item(6, 16).Value = 0.1486879047638899
$ws.Cells.Item(6, 17).Value = 153.8092335876867
$ws.Cells.Item(6, 18).Value = 1384.28310228918
$ws.Cells.Item(6, 19).Value = 0.0752625336028252
$ws.Cells.Item(6, 20).Value = 0.09479682863475326
$ws.Cells.Item(7, 9).Value = 0.03478077306145753
$ws.Cells.Item(7, 10).Value = 0.04034942645199305
$ws.Cells.Item(7, 13).Value = 38.10639333333333
$ws.Cells.Item(7, 14).Value = 114.31918
$ws.Cells.Item(7, 15).Value = 0.3831479157160237
$ws.Cells.Item(7, 16).Value = 0.4159903984418967
$ws.Cells.Item(7, 17).Value = 27.23386448533556
$ws.Cells.Item(7, 18).Value = 245.10478036802
$ws.Cells.Item(7, 19).Value = 0.01332618070548948
$ws.Cells.Item(7, 20).Value = 0.0167849739866666
$ws.Cells.Item(8, 9).Value = 0.03478077306145753
$ws.Cells.Item(8, 10).Value = 0.04034942645199305
$ws.Cells.Item(8, 15).Value = 0.09199521176963764
$ws.Cells.Item(8, 16).Value = 0.09988081163714851
$ws.Cells.Item(8, 19).Value = 0.003199664583300494
$ws.Cells.Item(8, 20).Value = 0.004030133463118496
$ws.Cells.Item(9, 9).Value = 0.03478077306145753
$ws.Cells.Item(9, 10).Value = 0.04034942645199305
$ws.Cells.Item(9, 13).Value = 15.023598
$ws.Cells.Item(9, 14).Value = 45.070794
$ws.Cells.Item(9, 15).Value = 0.1510575983904562
$ws.Cells.Item(9, 16).Value = 0.1640058785774412
$ws.Cells.Item(9, 17).Value = 10.737060010774
$ws.Cells.Item(9, 18).Value = 96.63354009696602
$ws.Cells.Item(9, 19).Value = 0.005253900048827249
$ws.Cells.Item(9, 20).Value = 0.006617543135354968
$ws.Cells.Item(10, 9).Value = 0.03478077306145753
$ws.Cells.Item(10, 10).Value = 0.04034942645199305
$ws.Cells.Item(10, 13).Value = 23.556204
$ws.Cells.Item(10, 14).Value = 47.112408
$ws.Cells.Item(10, 15).Value = 0.236850294013169
$ws.Cells.Item(10, 16).Value = 0.1714350065796238
$ws.Cells.Item(10, 17).Value = 16.835140022652
$ws.Cells.Item(10, 18).Value = 101.010840135912
$ws.Cells.Item(10, 19).Value = 0.008237836325611524
$ws.Cells.Item(10, 20).Value = 0.006917304189281478
$ws.Cells.Item(11, 9).Value = 0.03478077306145753
$ws.Cells.Item(11, 10).Value = 0.04034942645199305
$ws.Cells.Item(11, 13).Value = 13.62041
$ws.Cells.Item(11, 14).Value = 40.86123000000001
$ws.Cells.Item(11, 15).Value = 0.1369489801107134
$ws.Cells.Item(11, 16).Value = 0.1486879047638899
$ws.Cells.Item(11, 17).Value = 9.734230078663336
$ws.Cells.Item(11, 18).Value = 87.60807070797003
$ws.Cells.Item(11, 19).Value = 0.004763191398228783
$ws.Cells.Item(11, 20).Value = 0.005999471677571522
$ws.Cells.Item(12, 7).Value = 8.5075845
$ws.Cells.Item(12, 8).Value = 17.015169
$ws.Cells.Item(12, 9).Value = 0.4140321595768645
$ws.Cells.Item(12, 10).Value = 0.3202144691088791
$ws.Cells.Item(12, 13).Value = 38.10639333333333
$ws.Cells.Item(12, 14).Value = 114.31918
$ws.Cells.Item(12, 15).Value = 0.3831479157160237
$ws.Cells.Item(12, 16).Value = 0.4159903984418967
$ws.Cells.Item(12, 17).Value = 324.19336127357
$ws.Cells.Item(12, 18).Value = 1945.16016764142
$ws.Cells.Item(12, 19).Value = 0.1586355589812798
$ws.Cells.Item(12, 20).Value = 0.133206144591463
$ws.Cells.Item(13, 7).Value = 8.5075845
$ws.Cells.Item(13, 8).Value = 17.015169
$ws.Cells.Item(13, 9).Value = 0.4140321595768645
$ws.Cells.Item(13, 10).Value = 0.3202144691088791
$ws.Cells.Item(13, 15).Value = 0.09199521176963764
$ws.Cells.Item(13, 16).Value = 0.09988081163714851
$ws.Cells.Item(13, 17).Value = 77.840008261398
$ws.Cells.Item(13, 18).Value = 467.040049568388
$ws.Cells.Item(13, 19).Value = 0.03808897619971405
$ws.Cells.Item(13, 20).Value = 0.03198328107255347
$ws.Cells.Item(14, 7).Value = 8.5075845
$ws.Cells.Item(14, 8).Value = 17.015169
$ws.Cells.Item(14, 9).Value = 0.4140321595768645
$ws.Cells.Item(14, 10).Value = 0.3202144691088791
$ws.Cells.Item(14, 13).Value = 15.023598
$ws.Cells.Item(14, 14).Value = 45.070794
$ws.Cells.Item(14, 15).Value = 0.1510575983904562
$ws.Cells.Item(14, 16).Value = 0.1640058785774412
$ws.Cells.Item(14, 17).Value = 127.814529479031
$ws.Cells.Item(14, 18).Value = 766.887176874186
$ws.Cells.Item(14, 19).Value = 0.06254270368209526
$ws.Cells.Item(14, 20).Value = 0.05251705533941064
$ws.Cells.Item(15, 7).Value = 8.5075845
$ws.Cells.Item(15, 8).Value = 17.015169
$ws.Cells.Item(15, 9).Value = 0.4140321595768645
$ws.Cells.Item(15, 10).Value = 0.3202144691088791
$ws.Cells.Item(15, 13).Value = 23.556204
$ws.Cells.Item(15, 14).Value = 47.112408
$ws.Cells.Item(15, 15).Value = 0.236850294013169
$ws.Cells.Item(15, 16).Value = 0.1714350065796238
$ws.Cells.Item(15, 17).Value = 200.406396029238
$ws.Cells.Item(15, 18).Value = 801.625584116952
$ws.Cells.Item(15, 19).Value = 0.09806363872668766
$ws.Cells.Item(15, 20).Value = 0.05489596961857145
$ws.Cells.Item(16, 7).Value = 8.5075845
$ws.Cells.Item(16, 8).Value = 17.015169
$ws.Cells.Item(16, 9).Value = 0.4140321595768645
$ws.Cells.Item(16, 10).Value = 0.3202144691088791
$ws.Cells.Item(16, 13).Value = 13.62041
$ws.Cells.Item(16, 14).Value = 40.86123000000001
$ws.Cells.Item(16, 15).Value = 0.1369489801107134
$ws.Cells.Item(16, 16).Value = 0.1486879047638899
$ws.Cells.Item(16, 17).Value = 115.876788999645
$ws.Cells.Item(16, 18).Value = 695.2607339978701
$ws.Cells.Item(16, 19).Value = 0.05670128198708772
$ws.Cells.Item(16, 20).Value = 0.04761201848688058
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = 0.3333333333333333
$ws.Cells.Item(17, 7).Value = 0.03330533333333333
$ws.Cells.Item(17, 8).Value = 0.099916
$ws.Cells.Item(17, 9).Value = 0.001620845386305282
$ws.Cells.Item(17, 10).Value = 0.001880354458746944
$ws.Cells.Item(17, 13).Value = 38.10639333333333
$ws.Cells.Item(17, 14).Value = 114.31918
$ws.Cells.Item(17, 15).Value = 0.3831479157160237
$ws.Cells.Item(17, 16).Value = 0.4159903984418967
$ws.Cells.Item(17, 17).Value = 1.269146132097778
$ws.Cells.Item(17, 18).Value = 11.42231518888
$ws.Cells.Item(17, 19).Value = 0.0006210235314608022
$ws.Cells.Item(17, 20).Value = 0.000782209400506138
$ws.Cells.Item(18, 5).Value = 1
$ws.Cells.Item(18, 6).Value = 0.3333333333333333
$ws.Cells.Item(18, 7).Value = 0.03330533333333333
$ws.Cells.Item(18, 8).Value = 0.099916
$ws.Cells.Item(18, 9).Value = 0.001620845386305282
$ws.Cells.Item(18, 10).Value = 0.001880354458746944
$ws.Cells.Item(18, 15).Value = 0.09199521176963764
$ws.Cells.Item(18, 16).Value = 0.09988081163714851
$ws.Cells.Item(18, 17).Value = 0.304726614448
$ws.Cells.Item(18, 18).Value = 2.742539530032
$ws.Cells.Item(18, 19).Value = 0.0001491100145589946
$ws.Cells.Item(18, 20).Value = 0.0001878113295051758
$ws.Cells.Item(19, 5).Value = 1
$ws.Cells.Item(19, 6).Value = 0.3333333333333333
$ws.Cells.Item(19, 7).Value = 0.03330533333333333
$ws.Cells.Item(19, 8).Value = 0.099916
$ws.Cells.Item(19, 9).Value = 0.001620845386305282
$ws.Cells.Item(19, 10).Value = 0.001880354458746944
$ws.Cells.Item(19, 13).Value = 15.023598
$ws.Cells.Item(19, 14).Value = 45.070794
$ws.Cells.Item(19, 15).Value = 0.1510575983904562
$ws.Cells.Item(19, 16).Value = 0.1640058785774412
$ws.Cells.Item(19, 17).Value = 0.500365939256
$ws.Cells.Item(19, 18).Value = 4.503293453304
$ws.Cells.Item(19, 19).Value = 0.0002448410114175271
$ws.Cells.Item(19, 20).Value = 0.0003083891850438014
$ws.Cells.Item(20, 5).Value = 1
$ws.Cells.Item(20, 6).Value = 0.3333333333333333
$ws.Cells.Item(20, 7).Value = 0.03330533333333333
$ws.Cells.Item(20, 8).Value = 0.099916
$ws.Cells.Item(20, 9).Value = 0.001620845386305282
$ws.Cells.Item(20, 10).Value = 0.001880354458746944
$ws.Cells.Item(20, 13).Value = 23.556204
$ws.Cells.Item(20, 14).Value = 47.112408
$ws.Cells.Item(20, 15).Value = 0.236850294013169
$ws.Cells.Item(20, 16).Value = 0.1714350065796238
$ws.Cells.Item(20, 17).Value = 0.784547226288
$ws.Cells.Item(20, 18).Value = 4.707283357728
$ws.Cells.Item(20, 19).Value = 0.0003838977062962945
$ws.Cells.Item(20, 20).Value = 0.0003223585790073073
$ws.Cells.Item(21, 5).Value = 1
$ws.Cells.Item(21, 6).Value = 0.3333333333333333
$ws.Cells.Item(21, 7).Value = 0.03330533333333333
$ws.Cells.Item(21, 8).Value = 0.099916
$ws.Cells.Item(21, 9).Value = 0.001620845386305282
$ws.Cells.Item(21, 10).Value = 0.001880354458746944
$ws.Cells.Item(21, 13).Value = 13.62041
$ws.Cells.Item(21, 14).Value = 40.86123000000001
$ws.Cells.Item(21, 15).Value = 0.1369489801107134
$ws.Cells.Item(21, 16).Value = 0.1486879047638899
$ws.Cells.Item(21, 17).Value = 0.4536322951866667
$ws.Cells.Item(21, 18).Value = 4.082690656680001
$ws.Cells.Item(21, 19).Value = 0.0002219731225716636
$ws.Cells.Item(21, 20).Value = 0.0002795859646845212
